$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the Axonn 2-pin connector part number (C10 / E10) to the fixed one.
$ws.Range("C10").Value = "MSAP102SBS1G2GCEL2E"
$ws.Range("E10").Value = "MSAP102SBS1G2GCEL2E"

# The "qty for 5 boards" helper note (header + computed value) is no longer
# needed now that the part number has been corrected in place.
$ws.Range("I3").ClearContents()
$ws.Range("I10").ClearContents()

# The whole "corrected part number" helper column is now obsolete; remove it.
$ws.Columns("J").Delete()
